$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "quantity" header to "available"
$ws.Range("C1").Value = "available"

# Convert the quantity column into a boolean "available" column
$fmt = '"TRUE";"TRUE";"FALSE"'

$values = @(
    [bool]$true,
    [bool]$false,
    [bool]$true,
    [bool]$true,
    [bool]$true,
    [bool]$false,
    [bool]$false,
    [bool]$true,
    [bool]$true,
    [bool]$true
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $cell = $ws.Range("C$row")
    $cell.Value = $values[$i]
    $cell.NumberFormat = $fmt
}

$ws.Rows(1048576).RowHeight = 12.8

$ws.Range("C12").Select()
